# Insert a new row at row 335 (shifts existing rows 335:442 down to 336:443)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new record's data
$ws.Range("A335").Value = 6
$ws.Range("B335").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C335").Value = "Metropolitana"
$ws.Range("D335").Value = 44524
$ws.Range("E335").Value = 13
$ws.Range("F335").Value = 100112003
$ws.Range("G335").Value = "Ajo"
$ws.Range("H335").Value = "Chino"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 1900
$ws.Range("K335").Value = 16500
$ws.Range("L335").Value = 17000
$ws.Range("M335").Value = 16711
$ws.Range("N335").Value = "$/caja 10 kilos"
$ws.Range("O335").Value = "China"
$ws.Range("P335").Value = 1671
$ws.Range("Q335").Value = 10
$ws.Range("R335").Value = "Hortaliza"
